$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the date in A16 from 45926 (2025-09-26) to 45927 (2025-09-27)
$ws.Range("A16").Value = 45927

# Move the active selection to E15
$ws.Range("E15").Select()
